# repull data, push all data, mean calculation
# Updates the "dSF" (column F) values for a set of rows on Sheet1 to reflect
# repulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -1
    11 = 2
    12 = 1
    13 = -2
    14 = -1
    16 = 2
    17 = 2
    28 = -1
    31 = 0
    37 = -1
    40 = -1
    41 = 0
    45 = 10
    52 = -6
    55 = -4
    64 = 16
    68 = 6
    70 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
